$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 previously held "Dall Makhini" / 170 / a hyperlinked image URL.
# That row's dish data is removed, leaving only the (still hyperlink-styled)
# empty C5 cell; the now-unused shared strings drop out of the table and
# row 6 (butter naan dal makhani) shifts index references accordingly.
$ws.Range("A5:B5").ClearContents()
$ws.Range("C5").ClearContents()

# Update the active selection left behind by the edit.
$ws.Range("C13").Select()
